# Add a "months" list (Jan..Dec) to the hidden Dropdown sheet in column E,
# wire it up as a new data-validation list source for the Height Start/End
# Month columns (AU:AV) on the main "Report Input Template" sheet, and
# leave the selection state matching where the author ended up working.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Dropdown sheet: populate E1:E12 with the month names.
# ---------------------------------------------------------------------
$dropdown = $wb.Worksheets.Item("Dropdown")

$months = @("Jan", "Feb", "Mar", "Apr", "May", "Jun", "Jul", "Aug", "Sep", "Oct", "Nov", "Dec")
for ($i = 0; $i -lt $months.Length; $i++) {
    $dropdown.Cells.Item($i + 1, 5).Value = $months[$i]
}

# Matches the author's final selection on the Dropdown sheet.
[void]$dropdown.Range("E1:E12").Select()

# ---------------------------------------------------------------------
# 2. Report Input Template sheet: add the list validation for the new
#    Height Start Month / Height End Month columns (AU, AV), sourced
#    from the Dropdown sheet's new E1:E12 range.
# ---------------------------------------------------------------------
$main = $wb.Worksheets.Item("Report Input Template")

$rng = $main.Range("AU2:AV1048576")
$rng.Validation.Add(3, 1, 1, "=Dropdown!`$E`$1:`$E`$12")

# Matches the author's final selection on the main sheet.
[void]$main.Range("AR12").Select()
